$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 ---
$ws.Range("C2").Value = 2

# Apply the header style (bold/border/centered) to A2, and to the new A3:A5 cells too
$ws.Range("A1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)

# --- Row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "run-01"
$ws.Range("C3").Value = 5
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 1

# --- Row 4 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "run-01"
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = "--"
$ws.Range("J4").Value = "Test"

# --- Row 5 ---
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "run-01"
$ws.Range("C5").Value = 11
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = "--"
$ws.Range("F5").Value = "Interférence d'interaction"
$ws.Range("G5").Value = "Système (Machine)"
$ws.Range("H5").Value = "Test"
$ws.Range("I5").Value = "Neutre"
$ws.Range("J5").Value = "test`n"
$ws.Range("K5").Value = "Neutre"
$ws.Range("N5").Value = "Neutre"
